$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Locate the last populated row in column A (the date column) and append after it.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row   # xlUp
$firstNewRow = $lastRow + 1

# New daily data rows: date serial, nuovi pos., somma mobile 7gg., somma mobile 7gg. per 100mila ab.
# Covers 2021-05-28 through 2021-06-28 inclusive (aggiornamento fino a 28/06 incluso).
$newData = @(
    @(44344, 0, 0, 0),
    @(44345, 0, 0, 0),
    @(44346, 0, 0, 0),
    @(44347, 1, 1, 23.82654276864427),
    @(44348, 0, 1, 23.82654276864427),
    @(44349, 0, 1, 23.82654276864427),
    @(44350, 1, 2, 47.65308553728854),
    @(44351, 0, 2, 47.65308553728854),
    @(44352, 0, 2, 47.65308553728854),
    @(44353, 1, 3, 71.47962830593281),
    @(44354, 0, 2, 47.65308553728854),
    @(44355, 0, 2, 47.65308553728854),
    @(44356, 0, 2, 47.65308553728854),
    @(44357, 0, 1, 23.82654276864427),
    @(44358, 0, 1, 23.82654276864427),
    @(44359, 0, 1, 23.82654276864427),
    @(44360, 0, 0, 0),
    @(44361, 0, 0, 0),
    @(44362, 0, 0, 0),
    @(44363, 0, 0, 0),
    @(44364, 0, 0, 0),
    @(44365, 0, 0, 0),
    @(44366, 0, 0, 0),
    @(44367, 0, 0, 0),
    @(44368, 0, 0, 0),
    @(44369, 0, 0, 0),
    @(44370, 0, 0, 0),
    @(44371, 0, 0, 0),
    @(44372, 0, 0, 0),
    @(44373, 0, 0, 0),
    @(44374, 0, 0, 0),
    @(44375, 0, 0, 0)
)

$numRows = $newData.Count
$arr = New-Object "object[,]" $numRows,4
for ($i = 0; $i -lt $numRows; $i++) {
    for ($j = 0; $j -lt 4; $j++) {
        $arr[$i,$j] = $newData[$i][$j]
    }
}

$lastNewRow = $firstNewRow + $numRows - 1
$destRange = $ws.Range($ws.Cells.Item($firstNewRow,1), $ws.Cells.Item($lastNewRow,4))
$destRange.Value = $arr

# Copy formatting (style/number format/borders) from the last existing data row onto the new rows
$srcRange = $ws.Range($ws.Cells.Item($lastRow,1), $ws.Cells.Item($lastRow,4))
$srcRange.Copy()
$destRange.PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

"Added rows $firstNewRow to $lastNewRow"
